$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Existencias" (column E) values for rows 2-8
$ws.Range("E2").Value = 590.0
$ws.Range("E3").Value = 600.0
$ws.Range("E4").Value = 781.0
$ws.Range("E5").Value = 943.0
$ws.Range("E6").Value = 901.0
$ws.Range("E7").Value = 931.0
$ws.Range("E8").Value = 994.0

# Add new row 9: Vaca Lechera 2 / Planta del Prado
$ws.Range("A9").Value = "Vaca Lechera 2"
$ws.Range("B9").Value = "Planta del Prado"
$ws.Range("C9").Value = 145.0
$ws.Range("D9").Value = 3.0
$ws.Range("E9").Value = 600.0

# Add new row 10: Nutrimayma Fase 1 / Comayma
$ws.Range("A10").Value = "Nutrimayma Fase 1"
$ws.Range("B10").Value = "Comayma"
$ws.Range("C10").Value = 215.0
$ws.Range("D10").Value = 3.0
$ws.Range("E10").Value = 600.0
